# tippmix isolated from main file, moved to ML_PL_new
# Populate the "Team_tippmix" column (E) for the rows that were still
# missing it, reusing existing team/city names where they match and
# introducing the handful of genuinely new club names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E17").Value  = "Nottingham"
$ws.Range("E18").Value  = "Southampton"
$ws.Range("E90").Value  = "Nantes"
$ws.Range("E92").Value  = "Paris SG"
$ws.Range("E101").Value = "Boavista"
$ws.Range("E102").Value = "Casa Pia"
$ws.Range("E103").Value = "Estoril"
$ws.Range("E105").Value = "Famalicao"
$ws.Range("E107").Value = "Gil Vicente"
$ws.Range("E108").Value = "Guimaraes"
$ws.Range("E110").Value = "Nacional Madeira"
$ws.Range("E111").Value = "Porto"
$ws.Range("E113").Value = "Santa Clara"
$ws.Range("E115").Value = "Sporting CP"
$ws.Range("E118").Value = "Almere City"
$ws.Range("E119").Value = "Feyenoord"
$ws.Range("E120").Value = "F. Sittard"
$ws.Range("E121").Value = "Go Ahead Eagles"
$ws.Range("E123").Value = "Heerenveen"
$ws.Range("E125").Value = "NAC Breda"
$ws.Range("E126").Value = "Nijmegen"
$ws.Range("E127").Value = "PSV"
$ws.Range("E128").Value = "Sparta Rotterdam"
$ws.Range("E131").Value = "RKC Waalwijk"

# Widen column D to fit its longest entry now that the sheet has been
# reviewed end-to-end (matches the "bestFit" auto-sized column behaviour).
$ws.Columns.Item(4).ColumnWidth = 24.64

# Scroll the view back to the top and leave the selection on the first
# newly-filled cell instead of the far-down cell it was parked on before.
[void]$ws.Range("E19").Select()
